$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 161; this shifts the existing rows 161-201 down to 162-202
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row 161 with the new data record
$ws.Cells.Item(161, 1).Value = 11
$ws.Cells.Item(161, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(161, 3).Value = "Bíobío"
$ws.Cells.Item(161, 4).Value = 44798
$ws.Cells.Item(161, 5).Value = 8
$ws.Cells.Item(161, 6).Value = 100112003
$ws.Cells.Item(161, 7).Value = "Ajo"
$ws.Cells.Item(161, 8).Value = "Chino"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 400
$ws.Cells.Item(161, 11).Value = 21000
$ws.Cells.Item(161, 12).Value = 23000
$ws.Cells.Item(161, 13).Value = 22000
$ws.Cells.Item(161, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(161, 15).Value = "China"
$ws.Cells.Item(161, 16).Value = 2200
$ws.Cells.Item(161, 17).Value = 10
$ws.Cells.Item(161, 18).Value = "Hortaliza"
